$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Incorporate the new "Hommes 2019/20 JEL" reference as row 13, matching
# the layout of the existing rows (citation in column A, note in column C).
$ws.Range("A13").Value = "Hommes 2019/20 JEL"
$ws.Range("C13").Value = "oscillation under positive feedback"

# Column C entries use wrap text formatting throughout the sheet; apply the
# same formatting to the new note so it matches the rest of the table.
$ws.Range("C13").WrapText = $true

# Keep the selection on the newly added cell, consistent with the sheet's
# existing saved selection state.
$ws.Range("C13").Select() | Out-Null
